$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header row (row 2) relabelled from English to Spanish
$ws.Range("I2").Value = "Cargador"
$ws.Range("J2").Value = "Presentación"
$ws.Range("H2").Value = "Batería"
$ws.Range("G2").Value = "Características"
$ws.Range("F2").Value = "Potencia"
$ws.Range("E2").Value = "Voltaje"
$ws.Range("B2").Value = "Línea"

# Title cell: "CORDLESS MACHINES" -> "Power Tools"
$ws.Range("A1").Value = "Power Tools"

# Column I grew wider once its header text changed (no longer best-fit)
$ws.Columns.Item(9).ColumnWidth = 13.86

# Active selection moved from L3 to A2
$ws.Range("A2").Select() | Out-Null
